$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.137.23'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '3.306.23'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '528.12'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.72'
$ws.Range('E6').Value = '  -3.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.587'
$ws.Range('E7').Value = '  -2.85%  '
$ws.Range('D8').Value = '3.303.00'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.604'
$ws.Range('E10').Value = '  -2.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.30'
$ws.Range('E11').Value = '  -11.25%  '
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000258'
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.92'
$ws.Range('E14').Value = '  -2.40%  '
$ws.Range('D15').Value = '3.857.77'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.332.91'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.117'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '64.165.46'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.39'
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.16'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.953'
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '380.96'
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('E23').Value = '  +6.79%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.51'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('B25').Value = 'RenderToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.17'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.70'
$ws.Range('E26').Value = '  -3.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.09'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.72'
$ws.Range('E28').Value = '  +0.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.23'
$ws.Range('E29').Value = '  -3.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.17'
$ws.Range('E30').Value = '  -3.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.74'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '632.03'
$ws.Range('E32').Value = '  -3.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.65'
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.18'
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.105'
$ws.Range('E35').Value = '  -0.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.29'
$ws.Range('E36').Value = '  -3.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.30'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('E39').Value = '  -3.99%  '
$ws.Range('D40').Value = '0.0₃0751'
$ws.Range('E40').Value = '  +6.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.01'
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('E42').Value = '  +13.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.65'
$ws.Range('E43').Value = '  +6.23%  '
$ws.Range('E44').Value = '  -0.97%  '
$ws.Range('D45').Value = '2.945.61'
$ws.Range('E45').Value = '  +2.45%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.67'
$ws.Range('E46').Value = '  +0.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0399'
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.13'
$ws.Range('E48').Value = '  +2.50%  '
$ws.Range('E49').Value = '  -3.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '137.96'
$ws.Range('E50').Value = '  +2.16%  '
$ws.Range('E51').Value = '  -1.95%  '
